# Weekly refresh of the Fruta/Hortaliza (Frambuesa - Vega Monumental Concepción)
# price sheet: the data rows (2..21) get reshuffled among themselves -
# each destination row ends up with the D/L/M/N/O/P/S values that used to
# live in a different source row. All other columns (A,B,C,E..K,Q,R,T) are
# identical for every row already, so only these seven columns need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row data being shuffled.
$cols = @("D", "L", "M", "N", "O", "P", "S")

# Map: destination row -> source row (1-indexed data rows, 2..21).
$srcRow = @{
    2  = 7
    3  = 8
    4  = 15
    5  = 16
    6  = 9
    7  = 10
    8  = 11
    9  = 13
    10 = 14
    11 = 12
    12 = 20
    13 = 21
    14 = 17
    15 = 3
    16 = 4
    17 = 5
    18 = 6
    19 = 2
    20 = 18
    21 = 19
}

# Snapshot the original values for the columns involved before overwriting
# anything, so earlier writes don't clobber values still needed as sources.
# (Value2 is used for reads: this runtime's Value getter doesn't resolve.)
$orig = @{}
foreach ($r in 2..21) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

foreach ($r in 2..21) {
    $s = $srcRow[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $orig[$s][$c]
    }
}
